$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source dataset gained two additional weekly price records, which pushed
# all subsequent rows for this product/market block down (by 1 or 2 rows,
# depending on where the new records were spliced in) and grew the sheet
# from 109 to 111 data rows. Rather than replay each individual row-shift,
# insert two fresh rows at the bottom of the range to make room for the
# extra records and then rewrite the full A77:T111 block with its final
# contents in one pass.
$ws.Rows.Item(110).Insert()
$ws.Rows.Item(111).Insert()

$data = @(
    @(11, "Vega Monumental Concepción", "Bíobío", 44875, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Early Burlat", "Primera", 60, 24000, 25000, 24500, "`$/bandeja 10 kilos", "Provincia de Curicó", 2450, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44558, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Lapins", "Primera", 100, 4500, 5000, 4750, "`$/bandeja 10 kilos", "Región de O'Higgins", 475, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44558, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Santina", "Primera", 200, 4500, 5000, 4750, "`$/bandeja 10 kilos", "Región de O'Higgins", 475, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44537, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Lapins", "Primera", 160, 10000, 11000, 10500, "`$/caja 10 kilos", "Provincia de Curicó", 1050, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44537, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Royal Dawn", "Primera", 110, 14000, 15000, 14545, "`$/caja 10 kilos", "Provincia de Curicó", 1454, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44581, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Lapins", "Primera", 250, 5000, 5500, 5260, "`$/bandeja 10 kilos", "Provincia de Curicó", 526, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44526, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Santina", "Primera", 220, 15000, 16000, 15545, "`$/bandeja 5 kilos", "Provincia de Curicó", 3109, 5),
    @(11, "Vega Monumental Concepción", "Bíobío", 44203, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Lapins", "Primera", 100, 9000, 10000, 9500, "`$/caja 10 kilos", "Región de Ñuble", 950, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44203, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Lapins", "Segunda", 50, 8000, 8000, 8000, "`$/caja 10 kilos", "Región de Ñuble", 800, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44582, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Lapins", "Primera", 150, 6000, 6500, 6233, "`$/bandeja 10 kilos", "Provincia de Curicó", 623, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44187, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Lapins", "Primera", 200, 9000, 10000, 9500, "`$/caja 10 kilos", "Región de O'Higgins", 950, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44187, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Lapins", "Segunda", 100, 8000, 8000, 8000, "`$/caja 10 kilos", "Región de O'Higgins", 800, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44187, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Rainier", "Primera", 100, 10000, 11000, 10500, "`$/caja 10 kilos", "Región de O'Higgins", 1050, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44187, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Rainier", "Segunda", 50, 8000, 8000, 8000, "`$/caja 10 kilos", "Región de O'Higgins", 800, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44561, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Lapins", "Primera", 100, 5000, 6000, 5500, "`$/bandeja 10 kilos", "Región de O'Higgins", 550, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44561, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Lapins", "Segunda", 50, 4000, 4000, 4000, "`$/bandeja 10 kilos", "Región de O'Higgins", 400, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44166, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Early Burlat", "Primera", 200, 14000, 15000, 14500, "`$/caja 10 kilos", "Región de O'Higgins", 1450, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44166, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Early Burlat", "Segunda", 100, 12000, 12000, 12000, "`$/caja 10 kilos", "Región de O'Higgins", 1200, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44166, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Rainier", "Primera", 200, 15000, 16000, 15500, "`$/caja 10 kilos", "Región de O'Higgins", 1550, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44518, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Santina", "Segunda", 50, 28000, 28000, 28000, "`$/caja 10 kilos", "Provincia de Curicó", 2800, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44168, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Bing", "Primera", 200, 15000, 16000, 15500, "`$/caja 10 kilos", "Región de O'Higgins", 1550, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44168, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Rainier", "Segunda", 100, 13000, 13000, 13000, "`$/caja 10 kilos", "Región de O'Higgins", 1300, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44168, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Rainier", "Primera", 200, 15000, 16000, 15500, "`$/caja 10 kilos", "Región de O'Higgins", 1550, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44169, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Bing", "Primera", 200, 15000, 16000, 15500, "`$/caja 10 kilos", "Región de O'Higgins", 1550, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44169, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Bing", "Segunda", 100, 13000, 13000, 13000, "`$/caja 10 kilos", "Región de O'Higgins", 1300, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44169, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Rainier", "Primera", 100, 16000, 17000, 16500, "`$/caja 10 kilos", "Región de O'Higgins", 1650, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44169, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Rainier", "Segunda", 50, 14000, 14000, 14000, "`$/caja 10 kilos", "Región de O'Higgins", 1400, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44517, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Royal Dawn", "Segunda", 100, 20000, 22000, 21000, "`$/caja 10 kilos", "Provincia de Curicó", 2100, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44209, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Lapins", "Primera", 100, 9000, 10000, 9500, "`$/caja 10 kilos", "Región de Ñuble", 950, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44209, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Lapins", "Segunda", 80, 8000, 8000, 8000, "`$/caja 10 kilos", "Región de Ñuble", 800, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44579, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Lapins", "Especial", 170, 9000, 9500, 9235, "`$/caja 15 kilos", "Provincia de Curicó", 616, 15),
    @(11, "Vega Monumental Concepción", "Bíobío", 44579, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Santina", "Segunda", 180, 4000, 4500, 4278, "`$/bandeja 10 kilos", "Provincia de Curicó", 428, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44579, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Sweet Heart", "Primera", 250, 4500, 5000, 4700, "`$/bandeja 10 kilos", "Provincia de Curicó", 470, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44572, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Lapins", "Primera", 200, 5500, 6000, 5750, "`$/bandeja 10 kilos", "Región de Ñuble", 575, 10),
    @(11, "Vega Monumental Concepción", "Bíobío", 44572, 8, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Lapins", "Segunda", 100, 5000, 5000, 5000, "`$/bandeja 10 kilos", "Región de Ñuble", 500, 10)
)

$startRow = 77
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $targetRow = $startRow + $i
    for ($col = 1; $col -le $rowValues.Count; $col++) {
        $ws.Cells.Item($targetRow, $col).Value = $rowValues[$col - 1]
    }
}
